$wb = $excel.ActiveWorkbook

# --- Typography sheet: row 5 ("Typography_label") gains Wildcard
#     Characters (col G) and Wildcard Ranges (col I) values ---
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G5").Value = " !”`"#*%&()'`$+-@_, .:;?/~±×÷•º``´{}©£€^®¥_=[]¡¢|\¿><"
$wsTypo.Range("I5").Value = "a-z,A-Z,0-9"

# --- Translation sheet: two new rows for the "send amount of liters /
#     parse response" feature ---
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B42").Value = "SingleUseId37"
$wsTrans.Range("C42").Value = "Typographies_button"
$wsTrans.Range("D42").Value = "Center"
$wsTrans.Range("E42").Value = "LTR"
$wsTrans.Range("F42").Value = "Ok"

$wsTrans.Range("B43").Value = "SingleUseId38"
$wsTrans.Range("C43").Value = "Typography_label"
$wsTrans.Range("D43").Value = "Center"
$wsTrans.Range("E43").Value = "LTR"
$wsTrans.Range("F43").Value = "<value>"
